$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# --- Widen column C to fit the new "Java call" examples ---
$ws.Columns.Item(3).ColumnWidth = 20.8

# --- Row 8 (POWER): status IP -> N/A ---
$ws.Range("D8").Value = "N/A"

# --- Row 15 (MIN): add Java call example, status N/A -> Done ---
$ws.Range("C15").Value = "min(a1, a2, …, aN)"
$ws.Range("D15").Value = "Done"

# --- Row 18 (FLOOR): status N/A -> IP ---
$ws.Range("D18").Value = "IP"

# --- Row 19 (CEIL): status N/A -> IP ---
$ws.Range("D19").Value = "IP"

# --- Row 20 (LN): status IP -> N/A ---
$ws.Range("D20").Value = "N/A"

# --- Row 21 (LG): status IP -> N/A ---
$ws.Range("D21").Value = "N/A"

# --- Row 22 (LOG): status IP -> N/A ---
$ws.Range("D22").Value = "N/A"

# --- Row 23 (LOGN): status IP -> N/A ---
$ws.Range("D23").Value = "N/A"

# --- Row 24 (RNG): status N/A -> Done ---
$ws.Range("D24").Value = "Done"

# --- Row 25 (RNG_FLOAT): status N/A -> IP ---
$ws.Range("D25").Value = "IP"

# --- New row 44: ROUND function ---
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "ROUND"
$ws.Range("C44").Value = "round(a, accuracy)"
$ws.Range("D44").Value = "N/A"

# --- New row 45: MOD_POWER function ---
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "MOD_POWER"
$ws.Range("C45").Value = "modPow(base, exp, mod)"
$ws.Range("D45").Value = "IP"

# --- Update the sheet view: scroll down and select the newly added row ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A45:XFD45").Select() | Out-Null
